$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.172.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.087.69'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.39%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.650'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.30'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '60.33'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.368'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0737'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.24'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.380.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.829'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.085.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.988.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -2.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.34%  '
$ws.Range('E29').Value = '  -8.90%  '
$ws.Range('E30').Value = '  -6.80%  '
$ws.Range('E31').Value = '  +16.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0605'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0904'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.65%  '
$ws.Range('E38').Value = '  -4.24%  '
$ws.Range('E39').Value = '  -8.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.32'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0223'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0873'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.306.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.258.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.52%  '
